# "9th Stab - Cosmetic Changes"
# A new watch-date column is inserted before the existing "Jun_13" column:
# two brand-new columns (headers "Jun_17" and "Jun_15") are added, the old
# "Jun_13" header moves one slot to the right, and the old rating-history
# column (header "Jun_10") moves from C to E. The two freshly inserted
# columns (C, D) are backfilled with the default "UN" placeholder value,
# which is what column C used to hold for most rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 27

# --- Capture the existing column C contents/formatting before overwriting ---
$oldC = @{}
for ($r = 1; $r -le $lastRow; $r++) {
    $oldC[$r] = $ws.Cells.Item($r, 3).Value()
}
$colorIdxC10 = $ws.Range("C10").Interior.ColorIndex()
$colorIdxC17 = $ws.Range("C17").Interior.ColorIndex()

# --- Clear any inherited formatting from the two newly inserted columns ---
$ws.Range("C1:D27").ClearFormats()

# --- Move the old rating-history column (C) to its new home (E) ---
for ($r = 1; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 5).Value = $oldC[$r]
}
$ws.Range("E10").Interior.ColorIndex = $colorIdxC10
$ws.Range("E17").Interior.ColorIndex = $colorIdxC17

# --- Fill the two newly inserted columns (C, D) with the default "UN" value ---
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# --- Header row: shift "Jun_13" from B1 to D1, add the two new date headers ---
$ws.Range("D1").Value = "Jun_13"
$ws.Range("C1").Value = "Jun_15"
$ws.Range("B1").Value = "Jun_17"

# --- Match column widths of the old (narrow) watch-date column ---
$ws.Columns.Item(4).ColumnWidth = 7.14
$ws.Columns.Item(5).ColumnWidth = 7.14

# --- Cosmetically group the two superseded watch-date columns ---
$ws.Range("C1:D1").EntireColumn.Group()
